# "add tabel format baru" - renumber the table captions (4.2.3/4.2.4/4.2.5
# become 4.2.5/4.2.6/4.2.7) and roll the reporting year from 2020 to 2021
# across all four panels of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Panel 1 (columns A:M) ---------------------------------------------
$ws.Range("B1").Value = "Jumlah Tenaga Kesehatan Menurut Kelurahan/Desa in Kecamatan Watubangga. 2021"
$ws.Range("H1").Value = "Tabel 4.2.5"
$ws.Range("I1").Value = "Banyaknya Bayi yang Diimunisasi Menurut Jenis dan Desa/Kelurahan di Kecamatan Watubangga, 2021"
$ws.Range("B2").Value = "Number of Medical Personnel by Kelurahan/ Village in Watubangga Subdistrict, 2021"
$ws.Range("I2").Value = "Number of Immunized Babies by Types of Immunization and Kelurahan/Village Watubangga Subdistrict, 2021"

# --- Panel 2 (columns P:S), title keeps its "Tabel" / " 4.2.x." run split
$p1 = $ws.Range("P1")
$p1.Value = "Tabel 4.2.6."
$p1.Characters(6, 7).Font.Underline = $false
$ws.Range("Q1").Value = "Banyaknya Ibu Melahirkan dan Kelahiran Ditolong Tenaga Kesehatan Menurut Desa/Kelurahan di Kecamatan Watubangga, 2021"
$ws.Range("Q2").Value = "Number of Woman Giving Brth and Birth Assisted by Paramedics by Kelurahan/Village in Watubangga Subdistrict, 2021"

# --- Panel 3 (columns W:Z), title keeps its "Tabel" / " 4.2.x." run split
$w1 = $ws.Range("W1")
$w1.Value = "Tabel 4.2.7."
$w1.Characters(6, 7).Font.Underline = $false
$ws.Range("X1").Value = "Banyaknya Pasangan Usia Subur dan Peserta KB Menurut Desa/Kelurahan di Kecamatan Watubangga, 2021"
$ws.Range("X2").Value = "Number of Fertile Age Couples and Family Planning Members by Kelurahan/Village in Watubangga Subdistrict, 2021"

# --- Reset the saved view: scroll back to the top and drop the stray
#     mid-sheet selection that had been left over from editing.
$ws.Range("A1").Select()
